$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.639.95"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.037.04"
$ws.Range("E3").Value = "  +2.53%  "
$s_D4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $s_D4
$ws.Range("E4").Value = "  +0.02%  "
$s_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "383.87"
$ws.Range("D5").Style = $s_D5
$ws.Range("E5").Value = "  +0.86%  "
$s_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.77"
$ws.Range("D6").Style = $s_D6
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -0.03%  "
$s_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = $s_D9
$ws.Range("E9").Value = "  +0.27%  "
$s_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.74"
$ws.Range("D10").Style = $s_D10
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "3.513.11"
$ws.Range("E13").Value = "  +2.82%  "
$s_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.53"
$ws.Range("D14").Style = $s_D14
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "3.025.70"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$s_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.80"
$ws.Range("D17").Style = $s_D17
$ws.Range("E17").Value = "  -11.11%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$s_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.974"
$ws.Range("D18").Style = $s_D18
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("D19").Value = "51.653.60"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  -0.85%  "
$s_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("D21").Style = $s_D21
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.28%  "
$s_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.32"
$ws.Range("D24").Style = $s_D24
$ws.Range("E24").Value = "  -0.63%  "
$s_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.18"
$ws.Range("D25").Style = $s_D25
$ws.Range("E25").Value = "  -4.59%  "
$s_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.45"
$ws.Range("D26").Style = $s_D26
$ws.Range("E26").Value = "  +5.86%  "
$s_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.43"
$ws.Range("D27").Style = $s_D27
$ws.Range("E27").Value = "  +4.39%  "
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  +0.05%  "
$s_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.25"
$ws.Range("D30").Style = $s_D30
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("E33").Value = "  -2.71%  "
$s_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.05"
$ws.Range("D34").Style = $s_D34
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -1.47%  "
$s_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0447"
$ws.Range("D36").Style = $s_D36
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +7.90%  "
$s_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.07"
$ws.Range("D40").Style = $s_D40
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  +1.57%  "
$s_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").Style = $s_D42
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$s_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("D43").Style = $s_D43
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$s_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.86"
$ws.Range("D44").Style = $s_D44
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +3.84%  "
$s_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.66"
$ws.Range("D46").Style = $s_D46
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  +2.63%  "
$s_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("D48").Style = $s_D48
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "2.027.29"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "3.336.01"
$ws.Range("E50").Value = "  +2.55%  "
$s_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.516"
$ws.Range("D51").Style = $s_D51
$ws.Range("E51").Value = "  +6.34%  "
